# Apply weekly update to the Betarraga / Vega Modelo de Temuco sheet.
# The new weekly reading goes into row 192 (pushing previously-192..248
# readings down by one row, each keeping its own D/J-Q data), and the
# row that used to be last (248) lands in a brand new row 249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ D=44463; J=30; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44196; J=30; K=7000; L=7000; M=7000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=583; Q=12 },
    @{ D=44301; J=120; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44251; J=205; K=8000; L=9000; M=8244; N='$/docena de paquetes'; O='Provincia de Cautín'; P=687; Q=12 },
    @{ D=44243; J=125; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44252; J=155; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44166; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44166; J=70; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44168; J=185; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44369; J=90; K=7000; L=8000; M=7556; N='$/docena de paquetes'; O='Provincia de Cautín'; P=630; Q=12 },
    @{ D=44221; J=110; K=7000; L=7000; M=7000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=583; Q=12 },
    @{ D=44371; J=50; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44371; J=200; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44316; J=40; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44316; J=90; K=7000; L=8000; M=7556; N='$/docena de paquetes'; O='Región del Maule'; P=630; Q=12 },
    @{ D=44279; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44397; J=20; K=6000; L=6000; M=6000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=240; Q=25 },
    @{ D=44277; J=95; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44291; J=85; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44273; J=160; K=7000; L=8000; M=7500; N='$/docena de paquetes'; O='Provincia de Cautín'; P=625; Q=12 },
    @{ D=44438; J=50; K=9000; L=9000; M=9000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=750; Q=12 },
    @{ D=44438; J=80; K=6000; L=6000; M=6000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=240; Q=25 },
    @{ D=44372; J=50; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44372; J=50; K=7000; L=8000; M=7600; N='$/docena de paquetes'; O='Región del Maule'; P=633; Q=12 },
    @{ D=44286; J=80; K=7000; L=8000; M=7500; N='$/docena de paquetes'; O='Provincia de Cautín'; P=625; Q=12 },
    @{ D=44209; J=180; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44356; J=30; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44356; J=20; K=8000; L=8000; M=8000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=320; Q=25 },
    @{ D=44160; J=40; K=8000; L=8000; M=8000; N='$/saco 25 kilos'; O='Región del Maule'; P=320; Q=25 },
    @{ D=44351; J=125; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44351; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44365; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44365; J=185; K=8000; L=8000; M=8000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=320; Q=25 },
    @{ D=44306; J=35; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región de La Araucanía'; P=667; Q=12 },
    @{ D=44306; J=55; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44215; J=40; K=7000; L=7000; M=7000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=583; Q=12 },
    @{ D=44175; J=160; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44175; J=300; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44175; J=40; K=8000; L=8000; M=8000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=320; Q=25 },
    @{ D=44357; J=90; K=7000; L=8000; M=7556; N='$/docena de paquetes'; O='Provincia de Cautín'; P=630; Q=12 },
    @{ D=44357; J=50; K=8000; L=8000; M=8000; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=320; Q=25 },
    @{ D=44203; J=170; K=7000; L=8000; M=7471; N='$/docena de paquetes'; O='Provincia de Cautín'; P=623; Q=12 },
    @{ D=44162; J=60; K=7000; L=8000; M=7500; N='$/docena de paquetes'; O='Provincia de Cautín'; P=625; Q=12 },
    @{ D=44162; J=170; K=7000; L=8000; M=7471; N='$/docena de paquetes'; O='Región del Maule'; P=623; Q=12 },
    @{ D=44410; J=90; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44411; J=40; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44257; J=90; K=7000; L=8000; M=7556; N='$/docena de paquetes'; O='Provincia de Cautín'; P=630; Q=12 },
    @{ D=44244; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44176; J=80; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44176; J=200; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44239; J=210; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44376; J=65; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Región del Maule'; P=667; Q=12 },
    @{ D=44292; J=45; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44358; J=40; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44211; J=250; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 },
    @{ D=44425; J=30; K=9000; L=9000; M=9000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=750; Q=12 },
    @{ D=44425; J=40; K=8000; L=9000; M=8500; N='$/saco 25 kilos'; O='Provincia de Cautín'; P=340; Q=25 },
    @{ D=44323; J=55; K=8000; L=8000; M=8000; N='$/docena de paquetes'; O='Provincia de Cautín'; P=667; Q=12 }
)

$startRow = 192
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 4).Value = $data.D           # D: Fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 10).Value = $data.J           # J: Volumen
    $ws.Cells.Item($r, 11).Value = $data.K           # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $data.L           # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $data.M           # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $data.N           # N: Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $data.O           # O: Origen
    $ws.Cells.Item($r, 16).Value = $data.P           # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $data.Q           # Q: Kg o Unidades
}

# Row 249 is brand new - the columns that are constant across every row in
# this block (A, B, C, E, F, G, H, I, R) need to be (re)written explicitly
# since that row did not exist before.
$ws.Cells.Item(249, 1).Value = 10
$ws.Cells.Item(249, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(249, 3).Value = "La Araucanía"
$ws.Cells.Item(249, 5).Value = 9
$ws.Cells.Item(249, 6).Value = 100114014
$ws.Cells.Item(249, 7).Value = "Betarraga"
$ws.Cells.Item(249, 8).Value = "Sin especificar"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 18).Value = "Hortaliza"
